$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) with corrected values
# (Correccion a Diebold Mariano)
$ws.Range("C2").Value = -1.369222752465303
$ws.Range("D2").Value = 0.1847468298328279

$ws.Range("C3").Value = -0.9325905165252423
$ws.Range("D3").Value = 0.3611542257068452

$ws.Range("C4").Value = -0.5812698454246622
$ws.Range("D4").Value = 0.5669669067847858

$ws.Range("C5").Value = 0.08117548879069721
$ws.Range("D5").Value = 0.9360363648816095

$ws.Range("C6").Value = 0.3884104703610011
$ws.Range("D6").Value = 0.7014455791256897

$ws.Range("C7").Value = 0.8090514090100458
$ws.Range("D7").Value = 0.4271492026379053

$ws.Range("C8").Value = 1.727081942419339
$ws.Range("D8").Value = 0.09816893846958363

$ws.Range("C9").Value = 0.3346064155803629
$ws.Range("D9").Value = 0.7410911263585525

$ws.Range("C10").Value = 0.810420190048305
$ws.Range("D10").Value = 0.4263788909805717

$ws.Range("C11").Value = 0.674182658586801
$ws.Range("D11").Value = 0.507218003366757
